# Update the grade distribution counts on the "SECOND YEAR" sheet.
# These are the only real data values that changed between the two
# revisions of the workbook (the remaining instructor input cells
# across the other year sheets stayed blank/unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SECOND YEAR")

$ws.Range("E8").Value = 24
$ws.Range("E9").Value = 20
$ws.Range("E11").Value = 5
$ws.Range("E12").Value = 2
$ws.Range("E16").Value = 1
$ws.Range("E19").Value = 0

$excel.CalculateFull()
